$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: insert blank rows to make room for the new breakdown sub-rows ---
$ws.Rows("6:9").Insert()
$ws.Rows("11:14").Insert()
$ws.Rows("16:19").Insert()
$ws.Rows("21:24").Insert()

# --- Step 2: copy the bottom border / thick-bottom formatting from the old last row
#            (now row 20) down onto the new true last row (24) ---
$ws.Range("A20:E20").Copy() | Out-Null
$ws.Range("A24:E24").PasteSpecial(-4122) | Out-Null
$ws.Range("A20:E20").Borders.LineStyle = -4142
$excel.CutCopyMode = 0

# --- Step 3: write cell values ---
$ws.Cells.Item(5,1).Value = "Үй чарба жүргүзүү"
$ws.Cells.Item(5,2).Value = "Ведение домашнего хозяйства"
$ws.Cells.Item(5,3).Value = "Housekeeping"
$ws.Cells.Item(5,4).Value = 13.2
$ws.Cells.Item(5,5).Value = 11.5
$ws.Cells.Item(6,1).Value = "Шаар жерлери"
$ws.Cells.Item(6,2).Value = "Городские поселения"
$ws.Cells.Item(6,3).Value = "City"
$ws.Cells.Item(6,4).Value = 12.5
$ws.Cells.Item(6,5).Value = 10.7
$ws.Cells.Item(7,1).Value = "Айыл аймагы"
$ws.Cells.Item(7,2).Value = "Сельская местность"
$ws.Cells.Item(7,3).Value = "Village"
$ws.Cells.Item(7,4).Value = 13.9
$ws.Cells.Item(7,5).Value = 11.9
$ws.Cells.Item(8,1).Value = "Эркектер"
$ws.Cells.Item(8,2).Value = "Мужчины"
$ws.Cells.Item(8,3).Value = "Man"
$ws.Cells.Item(8,4).Value = 6.5
$ws.Cells.Item(8,5).Value = 3.9
$ws.Cells.Item(9,1).Value = "Аялдар"
$ws.Cells.Item(9,2).Value = "Женщины"
$ws.Cells.Item(9,3).Value = "Woman"
$ws.Cells.Item(9,4).Value = 18.8
$ws.Cells.Item(9,5).Value = 18.1
$ws.Cells.Item(10,1).Value = "Короо жанындагы, дача, бак участогуна карата жумуш"
$ws.Cells.Item(10,2).Value = "Работа на приусадебном, дачном, садовом участке"
$ws.Cells.Item(10,3).Value = "Work on a personal, country, garden plot"
$ws.Cells.Item(10,4).Value = 1.7
$ws.Cells.Item(10,5).Value = 2.6
$ws.Cells.Item(11,1).Value = "Шаар жерлери"
$ws.Cells.Item(11,2).Value = "Городские поселения"
$ws.Cells.Item(11,3).Value = "City"
$ws.Cells.Item(11,4).Value = 0.4
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(12,1).Value = "Айыл аймагы"
$ws.Cells.Item(12,2).Value = "Сельская местность"
$ws.Cells.Item(12,3).Value = "Village"
$ws.Cells.Item(12,4).Value = 3.1
$ws.Cells.Item(12,5).Value = 3.5
$ws.Cells.Item(13,1).Value = "Эркектер"
$ws.Cells.Item(13,2).Value = "Мужчины"
$ws.Cells.Item(13,3).Value = "Man"
$ws.Cells.Item(13,4).Value = 2.4
$ws.Cells.Item(13,5).Value = 3.6
$ws.Cells.Item(14,1).Value = "Аялдар"
$ws.Cells.Item(14,2).Value = "Женщины"
$ws.Cells.Item(14,3).Value = "Woman"
$ws.Cells.Item(14,4).Value = 1
$ws.Cells.Item(14,5).Value = 1.7
$ws.Cells.Item(15,1).Value = "Балдарды тарбиялоо"
$ws.Cells.Item(15,2).Value = "Воспитание детей"
$ws.Cells.Item(15,3).Value = "Parenting"
$ws.Cells.Item(15,4).Value = 0.9
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(16,1).Value = "Шаар жерлери"
$ws.Cells.Item(16,2).Value = "Городские поселения"
$ws.Cells.Item(16,3).Value = "City"
$ws.Cells.Item(16,4).Value = 0.9
$ws.Cells.Item(16,5).Value = 1.7
$ws.Cells.Item(17,1).Value = "Айыл аймагы"
$ws.Cells.Item(17,2).Value = "Сельская местность"
$ws.Cells.Item(17,3).Value = "Village"
$ws.Cells.Item(17,4).Value = 0.9
$ws.Cells.Item(17,5).Value = 2.2
$ws.Cells.Item(18,1).Value = "Эркектер"
$ws.Cells.Item(18,2).Value = "Мужчины"
$ws.Cells.Item(18,3).Value = "Man"
$ws.Cells.Item(18,4).Value = 0.6
$ws.Cells.Item(18,5).Value = 1.2
$ws.Cells.Item(19,1).Value = "Аялдар"
$ws.Cells.Item(19,2).Value = "Женщины"
$ws.Cells.Item(19,3).Value = "Woman"
$ws.Cells.Item(19,4).Value = 1.2
$ws.Cells.Item(19,5).Value = 2.8
$ws.Cells.Item(20,1).Value = "Туугандарга жана тааныштарга жардам берүү"
$ws.Cells.Item(20,2).Value = "Помощь родственникам и знакомым"
$ws.Cells.Item(20,3).Value = "Help for relatives and friends"
$ws.Cells.Item(20,4).Value = 0.5
$ws.Cells.Item(20,5).Value = 0.3
$ws.Cells.Item(21,1).Value = "Шаар жерлери"
$ws.Cells.Item(21,2).Value = "Городские поселения"
$ws.Cells.Item(21,3).Value = "City"
$ws.Cells.Item(21,4).Value = 0.4
$ws.Cells.Item(21,5).Value = 0.3
$ws.Cells.Item(22,1).Value = "Айыл аймагы"
$ws.Cells.Item(22,2).Value = "Сельская местность"
$ws.Cells.Item(22,3).Value = "Village"
$ws.Cells.Item(22,4).Value = 0.6
$ws.Cells.Item(22,5).Value = 0.3
$ws.Cells.Item(23,1).Value = "Эркектер"
$ws.Cells.Item(23,2).Value = "Мужчины"
$ws.Cells.Item(23,3).Value = "Man"
$ws.Cells.Item(23,4).Value = 0.6
$ws.Cells.Item(23,5).Value = 0.3
$ws.Cells.Item(24,1).Value = "Аялдар"
$ws.Cells.Item(24,2).Value = "Женщины"
$ws.Cells.Item(24,3).Value = "Woman"
$ws.Cells.Item(24,4).Value = 0.3
$ws.Cells.Item(24,5).Value = 0.3

# --- Step 4: formatting (bold category header rows) ---
$r = $ws.Range("A5:E5")
$r.Font.Bold = $true
$r.Font.Size = 9
$r.Font.Name = "Times New Roman"
$r = $ws.Range("A10:E10")
$r.Font.Bold = $true
$r.Font.Size = 9
$r.Font.Name = "Times New Roman"
$r = $ws.Range("A15:E15")
$r.Font.Bold = $true
$r.Font.Size = 9
$r.Font.Name = "Times New Roman"
$r = $ws.Range("A20:E20")
$r.Font.Bold = $true
$r.Font.Size = 9
$r.Font.Name = "Times New Roman"

# --- Step 5: wrap text ---
$ws.Range("A10:B10").WrapText = $true
$ws.Range("A11:E11").WrapText = $true
$ws.Range("A12:E12").WrapText = $true
$ws.Range("A13:E13").WrapText = $true
$ws.Range("A14:E14").WrapText = $true

# --- Step 6: explicit 0.0 number format on column E for rows where Excel shows it explicitly ---
$ws.Range("E15").NumberFormat = "0.0"
$ws.Range("E16").NumberFormat = "0.0"
$ws.Range("E17").NumberFormat = "0.0"
$ws.Range("E18").NumberFormat = "0.0"
$ws.Range("E19").NumberFormat = "0.0"

# --- Step 7: sheet view / dimension cosmetics ---
$ws.Range("A1").Select()

Write-Output "edit applied"
